# Apply the "Update countries & provincias Spain" data refresh.
# The underlying city list got re-sorted (several rows swap province
# names) and the per-province case counts (cols B-E) were refreshed,
# plus the "last updated" timestamp banner in A1 moved from 12:55 to 13:25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 13:25"

# Row 11: Alacant/Alicante -> Ciudad Real
$ws.Range("A11").Value = "Ciudad Real"
$ws.Range("B11").Value = 1755
$ws.Range("C11").Value = 236
$ws.Range("D11").Value = 1525
$ws.Range("E11").Value = 145

# Row 12: Ciudad Real -> Alacant/Alicante
$ws.Range("A12").Value = "Alacant/Alicante"
$ws.Range("B12").Value = 1734
$ws.Range("C12").Value = 45
$ws.Range("D12").Value = 1464
$ws.Range("E12").Value = 225

# Row 13: La Rioja -> Zaragoza
$ws.Range("A13").Value = "Zaragoza"
$ws.Range("B13").Value = 1449
$ws.Range("C13").Value = 101
$ws.Range("D13").Value = 1269
$ws.Range("E13").Value = 79

# Row 14: Zaragoza -> La Rioja
$ws.Range("A14").Value = "La Rioja"
$ws.Range("B14").Value = 1436
$ws.Range("C14").Value = 364
$ws.Range("D14").Value = 1007
$ws.Range("E14").Value = 65

# Row 15: Toledo -> Albacete
$ws.Range("A15").Value = "Albacete"
$ws.Range("B15").Value = 1386
$ws.Range("C15").Value = 236
$ws.Range("D15").Value = 1204
$ws.Range("E15").Value = 122

# Row 16: A Coruña -> A Coruña
$ws.Range("B16").Value = 1351
$ws.Range("C16").Value = 153
$ws.Range("D16").Value = 1261
$ws.Range("E16").Value = 40

# Row 17: Albacete -> Toledo
$ws.Range("A17").Value = "Toledo"
$ws.Range("B17").Value = 1317
$ws.Range("C17").Value = 236
$ws.Range("D17").Value = 1126
$ws.Range("E17").Value = 131

# Row 19: Malaga -> Pontevedra
$ws.Range("A19").Value = "Pontevedra"
$ws.Range("B19").Value = 1060
$ws.Range("C19").Value = 153
$ws.Range("D19").Value = 1005
$ws.Range("E19").Value = 9

# Row 20: Tenerife -> Malaga
$ws.Range("A20").Value = "Malaga"
$ws.Range("B20").Value = 1053
$ws.Range("C20").Value = 80
$ws.Range("D20").Value = 917
$ws.Range("E20").Value = 56

# Row 21: Cantabria -> Salamanca
$ws.Range("A21").Value = "Salamanca"
$ws.Range("B21").Value = 1030
$ws.Range("C21").Value = 157
$ws.Range("D21").Value = 774
$ws.Range("E21").Value = 99

# Row 22: Gipuzkoa/Guipuzcoa -> Tenerife
$ws.Range("A22").Value = "Tenerife"
$ws.Range("B22").Value = 1025
$ws.Range("C22").Value = 25
$ws.Range("D22").Value = 964
$ws.Range("E22").Value = 36

# Row 23: Caceres -> Cantabria
$ws.Range("A23").Value = "Cantabria"
$ws.Range("B23").Value = 1023
$ws.Range("C23").Value = 25
$ws.Range("D23").Value = 972
$ws.Range("E23").Value = 26

# Row 24: Pontevedra -> Gipuzkoa/Guipuzcoa
$ws.Range("A24").Value = "Gipuzkoa/Guipuzcoa"
$ws.Range("B24").Value = 1017
$ws.Range("C24").Value = 1503
$ws.Range("D24").Value = 630
$ws.Range("E24").Value = 34

# Row 25: Aragon -> Caceres
$ws.Range("A25").Value = "Caceres"
$ws.Range("B25").Value = 991
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 896
$ws.Range("E25").Value = 84

# Row 26: Salamanca -> Aragon
$ws.Range("A26").Value = "Aragon"
$ws.Range("B26").Value = 907
$ws.Range("C26").Value = 29
$ws.Range("D26").Value = 838
$ws.Range("E26").Value = 40

# Row 27: Murcia -> Valladolid
$ws.Range("A27").Value = "Valladolid"
$ws.Range("B27").Value = 886
$ws.Range("C27").Value = 127
$ws.Range("D27").Value = 702
$ws.Range("E27").Value = 57

# Row 28: Sevilla -> Murcia
$ws.Range("A28").Value = "Murcia"
$ws.Range("B28").Value = 872
$ws.Range("C28").Value = 16
$ws.Range("D28").Value = 836
$ws.Range("E28").Value = 20

# Row 29: Valladolid -> Sevilla
$ws.Range("A29").Value = "Sevilla"
$ws.Range("B29").Value = 830
$ws.Range("C29").Value = 13
$ws.Range("D29").Value = 791
$ws.Range("E29").Value = 26

# Row 30: Granada -> Leon
$ws.Range("A30").Value = "Leon"
$ws.Range("B30").Value = 821
$ws.Range("C30").Value = 118
$ws.Range("D30").Value = 626
$ws.Range("E30").Value = 77

# Row 31: Leon -> Granada
$ws.Range("A31").Value = "Granada"
$ws.Range("B31").Value = 806
$ws.Range("C31").Value = 11
$ws.Range("D31").Value = 746
$ws.Range("E31").Value = 49

# Row 32: Burgos -> Burgos
$ws.Range("B32").Value = 719
$ws.Range("C32").Value = 137
$ws.Range("D32").Value = 533
$ws.Range("E32").Value = 49

# Row 33: Castello/Castellon -> Segovia
$ws.Range("A33").Value = "Segovia"
$ws.Range("B33").Value = 567
$ws.Range("C33").Value = 131
$ws.Range("D33").Value = 378
$ws.Range("E33").Value = 58

# Row 34: Segovia -> Castello/Castellon
$ws.Range("A34").Value = "Castello/Castellon"
$ws.Range("B34").Value = 542
$ws.Range("C34").Value = 8
$ws.Range("D34").Value = 509
$ws.Range("E34").Value = 25

# Row 35: Jaen -> Guadalajara
$ws.Range("A35").Value = "Guadalajara"
$ws.Range("B35").Value = 535
$ws.Range("C35").Value = 236
$ws.Range("D35").Value = 436
$ws.Range("E35").Value = 86

# Row 36: Badajoz -> Soria
$ws.Range("A36").Value = "Soria"
$ws.Range("B36").Value = 523
$ws.Range("C36").Value = 61
$ws.Range("D36").Value = 432
$ws.Range("E36").Value = 30

# Row 37: Guadalajara -> Jaen
$ws.Range("A37").Value = "Jaen"
$ws.Range("B37").Value = 465
$ws.Range("C37").Value = 15
$ws.Range("D37").Value = 434
$ws.Range("E37").Value = 16

# Row 38: Soria -> Badajoz
$ws.Range("A38").Value = "Badajoz"
$ws.Range("B38").Value = 465
$ws.Range("C38").Value = 40
$ws.Range("D38").Value = 409
$ws.Range("E38").Value = 16

# Row 39: Cordoba -> Ourense
$ws.Range("A39").Value = "Ourense"
$ws.Range("B39").Value = 458
$ws.Range("C39").Value = 153
$ws.Range("D39").Value = 415
$ws.Range("E39").Value = 8

# Row 40: Cadiz -> Cordoba
$ws.Range("A40").Value = "Cordoba"
$ws.Range("B40").Value = 424
$ws.Range("C40").Value = 4
$ws.Range("D40").Value = 411
$ws.Range("E40").Value = 9

# Row 41: Ourense -> Avila
$ws.Range("A41").Value = "Avila"
$ws.Range("B41").Value = 414
$ws.Range("C41").Value = 82
$ws.Range("D41").Value = 291
$ws.Range("E41").Value = 41

# Row 42: Avila -> Cadiz
$ws.Range("A42").Value = "Cadiz"
$ws.Range("B42").Value = 406
$ws.Range("C42").Value = 8
$ws.Range("D42").Value = 391
$ws.Range("E42").Value = 7

# Row 43: Lugo -> Lugo
$ws.Range("B43").Value = 270
$ws.Range("C43").Value = 153
$ws.Range("D43").Value = 244
$ws.Range("E43").Value = 4

# Row 44: Gran Canaria -> Palencia
$ws.Range("A44").Value = "Palencia"
$ws.Range("B44").Value = 262
$ws.Range("C44").Value = 28
$ws.Range("D44").Value = 221
$ws.Range("E44").Value = 13

# Row 45: Cuenca -> Cuenca
$ws.Range("B45").Value = 253
$ws.Range("C45").Value = 236
$ws.Range("D45").Value = 180
$ws.Range("E45").Value = 55

# Row 46: Palencia -> Gran Canaria
$ws.Range("A46").Value = "Gran Canaria"
$ws.Range("B46").Value = 235
$ws.Range("C46").Value = 25
$ws.Range("D46").Value = 964
$ws.Range("E46").Value = 9

# Row 48: Teruel -> Teruel
$ws.Range("B48").Value = 208
$ws.Range("C48").Value = 10
$ws.Range("D48").Value = 186
$ws.Range("E48").Value = 12

# Row 49: Huesca -> Huesca
$ws.Range("B49").Value = 201
$ws.Range("C49").Value = 14
$ws.Range("D49").Value = 176
$ws.Range("E49").Value = 11

# Row 50: Almeria -> Zamora
$ws.Range("A50").Value = "Zamora"
$ws.Range("B50").Value = 192
$ws.Range("C50").Value = 30
$ws.Range("D50").Value = 144
$ws.Range("E50").Value = 18

# Row 51: Zamora -> Almeria
$ws.Range("A51").Value = "Almeria"
$ws.Range("B51").Value = 173
$ws.Range("C51").Value = 6
$ws.Range("D51").Value = 157
$ws.Range("E51").Value = 10

